$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Options" intent rows at the bottom of the intents table
$newRows = @(
    @("Options", "What can I do in this chatbot?"),
    @("Options", "What options do I hace?"),
    @("Options", "Tell me what can I ask for"),
    @("Options", "Show me what you got")
)

$startRow = 111
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Move/update the view the same way Excel would after scrolling to and
# selecting the new first empty row following the appended data
$lastRow = $startRow + $newRows.Length
$excel.ActiveWindow.ScrollRow = 94
$ws.Range("A" + $lastRow).Select()
